# Weekly data refresh: a new week's price record is inserted at the top of
# the data block (row 412), pushing all subsequent rows down by one and
# dropping the oldest historical rows off the bottom (old row 431 becomes
# the new row 432).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 412 - this shifts rows 412..431 down to 413..432
# and carries formatting (e.g. the date style on column D) down with them.
$ws.Rows.Item(412).Insert()

# Populate the newly inserted row 412 with the latest week's record.
$ws.Cells.Item(412, 1).Value = 6
$ws.Cells.Item(412, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(412, 3).Value = "Metropolitana"
$ws.Cells.Item(412, 4).Value = 45147
$ws.Cells.Item(412, 5).Value = 13
$ws.Cells.Item(412, 6).Value = 100112026
$ws.Cells.Item(412, 7).Value = "Haba"
$ws.Cells.Item(412, 8).Value = "Sin especificar"
$ws.Cells.Item(412, 9).Value = "Primera"
$ws.Cells.Item(412, 10).Value = 710
$ws.Cells.Item(412, 11).Value = 10000
$ws.Cells.Item(412, 12).Value = 12000
$ws.Cells.Item(412, 13).Value = 11014
$ws.Cells.Item(412, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(412, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(412, 16).Value = 441
$ws.Cells.Item(412, 17).Value = 25
$ws.Cells.Item(412, 18).Value = "Hortaliza"
